$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "2,3,4,5"
$ws.Range("D3").Value = "2,3,4,5"
$ws.Range("D4").Value = "3,5,4,2"
$ws.Range("C6").Value = $null
$ws.Range("D8").Value = "4,5"
$ws.Range("C9").Value = $null
$ws.Range("D9").Value = "5,7"

$ws.Range("C6").Select()
